# Update the dSF column (F) values on the active sheet to reflect
# re-pulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = 0
$ws.Range("F6").Value  = -7
$ws.Range("F7").Value  = -7
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 3
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 1
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("F29").Value = -1
